# Split the bold run "Process Specialist, Infosys Ltd." into three runs:
#   "Process Specialist, Infosys " + "BPM " + "Ltd."
# (all three keep identical Bold formatting), by typing "BPM " in front of
# "Ltd." and then toggling Bold off/on on the newly inserted text. The
# off/on toggle forces the engine to keep the freshly typed text as its own
# run instead of silently re-merging it with the neighboring identically
# formatted runs.

$d = $word.ActiveDocument

# 1. Find the word "Ltd." (with the trailing period) inside the job-title
#    line and collapse the found range to its start, i.e. the insertion
#    point immediately before "Ltd.".
$target = $d.Content
$found = $target.Find.Execute("Ltd.", $false, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Ltd.' in the document"
}
$target.Collapse(1)

# 2. Type the new text "BPM " right before "Ltd.". Word inserts it using
#    the formatting already in effect at that point (Bold/BCs), so the
#    paragraph text becomes "Process Specialist, Infosys BPM Ltd.".
$target.InsertBefore("BPM ")

# 3. Re-locate the just-inserted "BPM " text and briefly toggle Bold off
#    and back on. This keeps it as a distinct run (matching the target
#    OOXML, which stores it as its own <w:r>) while leaving its final
#    formatting (Bold) identical to its neighboring runs.
$bpmRange = $d.Content
$bpmRange.Find.Execute("BPM ", $false, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0) | Out-Null
$bpmRange.Bold = 0

$bpmRange2 = $d.Content
$bpmRange2.Find.Execute("BPM ", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null
$bpmRange2.Bold = 1
